$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 29.22839533333334
$ws.Range("N2").Value = 87.685186
$ws.Range("O2").Value = 0.4452121315669252
$ws.Range("P2").Value = 0.4988525987881078
$ws.Range("Q2").Value = 11.35234771866045
$ws.Range("R2").Value = 102.171129467944
$ws.Range("S2").Value = 0.4452121315669252
$ws.Range("T2").Value = 0.4988525987881078

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.232204800700274
$ws.Range("P3").Value = 0.2601815181287206
$ws.Range("Q3").Value = 5.920929490879113
$ws.Range("R3").Value = 53.28836541791201
$ws.Range("S3").Value = 0.232204800700274
$ws.Range("T3").Value = 0.2601815181287206

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 21.177737
$ws.Range("N4").Value = 42.355474
$ws.Range("O4").Value = 0.3225830677328006
$ws.Range("P4").Value = 0.2409658830831714
$ws.Range("Q4").Value = 8.225461287782668
$ws.Range("R4").Value = 49.35276772669601
$ws.Range("S4").Value = 0.3225830677328006
$ws.Range("T4").Value = 0.2409658830831714
